# "Add initial files for VGA Top" — append a new time-record entry (row 20)
# for task "VGA Top" / "Create necessary files", mirroring the existing rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Start row 20 as a copy of row 19 so it inherits the same number formats /
# alignment (date style, time style, formula style, etc.) as the rest of the
# table before we overwrite the per-row data.
$ws.Range("A19:F19").Copy($ws.Range("A20:F20"))

# A20 already equals A19 ("1.4.2020") after the copy, so it is left as-is.

# New entry: 10:30 - 10:45, task "VGA Top", notes "Create necessary files".
$ws.Range("B20").Value2 = 0.4375
$ws.Range("C20").Value2 = 0.44791666666666669
$ws.Range("D20").Formula = "=C20-B20"
$ws.Range("E20").Value2 = "VGA Top"
$ws.Range("F20").Value2 = "Create necessary files"

# Match the author's last selection in the saved file.
$ws.Range("C21").Select()

$wb.Save()
